$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1 ("Sheet1"): only the selected cell in the sheet view changes ---
$ws1.Range("D23").Select()

# --- Sheet2 ("Sheet2") ---
# selection moves
$ws2.Range("F13").Select()
# column E width changes (no longer "best fit", now an explicit custom width ~16.44)
$ws2.Columns(5).ColumnWidth = 15.714285714285714
# E7 switches from a text/shared-string value to a real date serial value
$ws2.Range("E7").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("E7").Value = 36729

# --- Sheet3 ("Sheet3") ---
# selection moves
$ws3.Range("F12").Select()
# column E width changes (no longer "best fit", now an explicit custom width ~16.22)
$ws3.Columns(5).ColumnWidth = 15.571428571428571
# E8 switches from a text/shared-string value to a real date serial value
$ws3.Range("E8").NumberFormat = "yyyy\-mm\-dd;@"
$ws3.Range("E8").Value = 36640

# Restore Sheet1 as the active sheet/tab (it was active before the edits,
# and should remain so -- only the in-sheet selections of Sheet2/Sheet3
# should change, not which sheet tab is active).
$ws1.Select()
